# Applies the coin price/volume/hour table updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{Ref='D2'; Val='307.87'},
    @{Ref='E2'; Val='-1.37%'},
    @{Ref='G2'; Val='18'},
    @{Ref='D3'; Val='35.87'},
    @{Ref='E3'; Val='-5.08%'},
    @{Ref='G3'; Val='18'},
    @{Ref='D4'; Val='5.118'},
    @{Ref='E4'; Val='-0.70%'},
    @{Ref='G4'; Val='18'},
    @{Ref='D5'; Val='0.07685'},
    @{Ref='E5'; Val='-2.81%'},
    @{Ref='G5'; Val='18'},
    @{Ref='B6'; Val='KuCoinToken'},
    @{Ref='C6'; Val='https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'},
    @{Ref='D6'; Val='8.296'},
    @{Ref='E6'; Val='0.32%'},
    @{Ref='G6'; Val='18'},
    @{Ref='B7'; Val='FTXToken'},
    @{Ref='C7'; Val='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'},
    @{Ref='D7'; Val='1.843'},
    @{Ref='E7'; Val='-3.26%'},
    @{Ref='G7'; Val='18'},
    @{Ref='B8'; Val='BTSEToken'},
    @{Ref='C8'; Val='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'},
    @{Ref='D8'; Val='2.955'},
    @{Ref='E8'; Val='-4.62%'},
    @{Ref='G8'; Val='18'},
    @{Ref='B9'; Val='MXToken'},
    @{Ref='C9'; Val='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Ref='D9'; Val='0.9208'},
    @{Ref='E9'; Val='0.04%'},
    @{Ref='G9'; Val='18'},
    @{Ref='B10'; Val='LiechtensteinCryptoassetsExchange'},
    @{Ref='C10'; Val='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'},
    @{Ref='D10'; Val='0.1093'},
    @{Ref='E10'; Val='-10.26%'},
    @{Ref='G10'; Val='18'},
    @{Ref='B11'; Val='WazirX'},
    @{Ref='C11'; Val='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'},
    @{Ref='D11'; Val='0.1842'},
    @{Ref='E11'; Val='-4.63%'},
    @{Ref='G11'; Val='18'},
    @{Ref='B12'; Val='MandalaExchangeToken'},
    @{Ref='C12'; Val='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'},
    @{Ref='D12'; Val='0.08722'},
    @{Ref='E12'; Val='-4.42%'},
    @{Ref='G12'; Val='18'},
    @{Ref='B13'; Val='BitrueCoin'},
    @{Ref='C13'; Val='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'},
    @{Ref='D13'; Val='0.03333'},
    @{Ref='E13'; Val='-0.06%'},
    @{Ref='G13'; Val='18'},
    @{Ref='B14'; Val='BitMartToken'},
    @{Ref='C14'; Val='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'},
    @{Ref='D14'; Val='0.09515'},
    @{Ref='E14'; Val='-1.06%'},
    @{Ref='G14'; Val='18'},
    @{Ref='B15'; Val='BitForexToken'},
    @{Ref='C15'; Val='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'},
    @{Ref='D15'; Val='0.001384'},
    @{Ref='E15'; Val='0.49%'},
    @{Ref='G15'; Val='18'},
    @{Ref='B16'; Val='TigerCash'},
    @{Ref='C16'; Val='https://coinranking.com/coin/6hIn06L2+tigercash-tch'},
    @{Ref='D16'; Val='0.006141'},
    @{Ref='E16'; Val='5.64%'},
    @{Ref='G16'; Val='18'},
    @{Ref='B17'; Val='LEO'},
    @{Ref='C17'; Val='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{Ref='D17'; Val='3.363'},
    @{Ref='E17'; Val='-4.50%'},
    @{Ref='G17'; Val='18'},
    @{Ref='B18'; Val='GateToken'},
    @{Ref='C18'; Val='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'},
    @{Ref='D18'; Val='4.400'},
    @{Ref='E18'; Val='-0.29%'},
    @{Ref='G18'; Val='18'},
    @{Ref='E19'; Val='-0.23%'},
    @{Ref='G19'; Val='18'},
    @{Ref='D20'; Val='6.283'},
    @{Ref='E20'; Val='18.98%'},
    @{Ref='G20'; Val='18'},
    @{Ref='D21'; Val='0.1292'},
    @{Ref='E21'; Val='1.56%'},
    @{Ref='G21'; Val='18'},
    @{Ref='D22'; Val='0.2317'},
    @{Ref='E22'; Val='-10.51%'},
    @{Ref='G22'; Val='18'},
    @{Ref='D23'; Val='0.04337'},
    @{Ref='E23'; Val='-0.56%'},
    @{Ref='G23'; Val='18'},
    @{Ref='D24'; Val='0.001202'},
    @{Ref='E24'; Val='-3.60%'},
    @{Ref='G24'; Val='18'},
    @{Ref='D25'; Val='0.004244'},
    @{Ref='G25'; Val='18'},
    @{Ref='E26'; Val='9.13%'},
    @{Ref='G26'; Val='18'},
    @{Ref='D27'; Val='0.0002908'},
    @{Ref='G27'; Val='18'},
    @{Ref='G28'; Val='18'},
    @{Ref='G29'; Val='18'},
    @{Ref='G30'; Val='18'},
    @{Ref='G31'; Val='18'},
    @{Ref='G32'; Val='18'},
    @{Ref='G33'; Val='18'},
    @{Ref='G34'; Val='18'},
    @{Ref='G35'; Val='18'},
    @{Ref='G36'; Val='18'},
    @{Ref='G37'; Val='18'},
    @{Ref='G38'; Val='18'},
    @{Ref='D39'; Val='0.02078'},
    @{Ref='E39'; Val='-3.12%'},
    @{Ref='G39'; Val='18'},
    @{Ref='D40'; Val='0.04917'},
    @{Ref='E40'; Val='-4.67%'},
    @{Ref='G40'; Val='18'},
    @{Ref='D41'; Val='0.007519'},
    @{Ref='E41'; Val='-0.48%'},
    @{Ref='G41'; Val='18'},
    @{Ref='D42'; Val='0.1347'},
    @{Ref='E42'; Val='-1.01%'},
    @{Ref='G42'; Val='18'},
    @{Ref='D43'; Val='0.008484'},
    @{Ref='E43'; Val='-6.70%'},
    @{Ref='G43'; Val='18'},
    @{Ref='E44'; Val='5.83%'},
    @{Ref='G44'; Val='18'},
    @{Ref='D45'; Val='0.008416'},
    @{Ref='E45'; Val='-2.47%'},
    @{Ref='G45'; Val='18'},
    @{Ref='D46'; Val='0.00006306'},
    @{Ref='E46'; Val='-5.83%'},
    @{Ref='G46'; Val='18'},
    @{Ref='D47'; Val='0.00000000752'},
    @{Ref='E47'; Val='0.33%'},
    @{Ref='G47'; Val='18'},
    @{Ref='D48'; Val='0.002857'},
    @{Ref='E48'; Val='-13.28%'},
    @{Ref='G48'; Val='18'},
    @{Ref='D49'; Val='0.001444'},
    @{Ref='E49'; Val='20.40%'},
    @{Ref='G49'; Val='18'},
    @{Ref='D50'; Val='0.00002105'},
    @{Ref='E50'; Val='0.33%'},
    @{Ref='G50'; Val='18'},
    @{Ref='D51'; Val='0.0002005'},
    @{Ref='E51'; Val='0.33%'},
    @{Ref='G51'; Val='18'}
)

foreach ($e in $edits) {
    $r = $ws.Range($e.Ref)
    $r.NumberFormat = '@'
    $r.Value = $e.Val
    $r.Style = 'Normal'
}
